{"js": "// Replace the multiplication-problem text in each table cell with the\n// new values from the commit's diff. Every \"before\" string is unique in\n// the document, so a scoped search-and-replace per pair is safe.\nconst replacements = [\n  [\"43\u00d782=\", \"96\u00d766=\"],\n  [\"38\u00d728=\", \"20\u00d734=\"],\n  [\"76\u00d755=\", \"71\u00d775=\"],\n  [\"43\u00d737=\", \"44\u00d763=\"],\n  [\"61\u00d736=\", \"32\u00d780=\"],\n  [\"60\u00d775=\", \"14\u00d776=\"],\n  [\"62\u00d763=\", \"20\u00d723=\"],\n  [\"53\u00d722=\", \"13\u00d760=\"],\n  [\"86\u00d718=\", \"94\u00d722=\"],\n  [\"30\u00d731=\", \"71\u00d745=\"],\n  [\"74\u00d768=\", \"56\u00d757=\"],\n  [\"30\u00d784=\", \"11\u00d787=\"],\n  [\"73\u00d760=\", \"91\u00d761=\"],\n  [\"80\u00d731=\", \"98\u00d754=\"],\n  [\"26\u00d753=\", \"71\u00d741=\"],\n  [\"18\u00d779=\", \"22\u00d798=\"],\n  [\"53\u00d718=\", \"86\u00d771=\"],\n  [\"81\u00d786=\", \"90\u00d776=\"],\n  [\"93\u00d763=\", \"52\u00d711=\"],\n  [\"15\u00d761=\", \"25\u00d725=\"],\n  [\"44\u00d751=\", \"44\u00d740=\"],\n  [\"47\u00d766=\", \"68\u00d723=\"],\n  [\"79\u00d783=\", \"83\u00d731=\"],\n  [\"50\u00d769=\", \"82\u00d725=\"],\n  [\"52\u00d759=\", \"40\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of replacements) {\n  const found = body.search(before, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the\n# new values from the commit's diff. Every \"before\" string is unique in\n# the document, so a single Find/Replace pass per pair is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"43\u00d782=\", \"96\u00d766=\"),\n    @(\"38\u00d728=\", \"20\u00d734=\"),\n    @(\"76\u00d755=\", \"71\u00d775=\"),\n    @(\"43\u00d737=\", \"44\u00d763=\"),\n    @(\"61\u00d736=\", \"32\u00d780=\"),\n    @(\"60\u00d775=\", \"14\u00d776=\"),\n    @(\"62\u00d763=\", \"20\u00d723=\"),\n    @(\"53\u00d722=\", \"13\u00d760=\"),\n    @(\"86\u00d718=\", \"94\u00d722=\"),\n    @(\"30\u00d731=\", \"71\u00d745=\"),\n    @(\"74\u00d768=\", \"56\u00d757=\"),\n    @(\"30\u00d784=\", \"11\u00d787=\"),\n    @(\"73\u00d760=\", \"91\u00d761=\"),\n    @(\"80\u00d731=\", \"98\u00d754=\"),\n    @(\"26\u00d753=\", \"71\u00d741=\"),\n    @(\"18\u00d779=\", \"22\u00d798=\"),\n    @(\"53\u00d718=\", \"86\u00d771=\"),\n    @(\"81\u00d786=\", \"90\u00d776=\"),\n    @(\"93\u00d763=\", \"52\u00d711=\"),\n    @(\"15\u00d761=\", \"25\u00d725=\"),\n    @(\"44\u00d751=\", \"44\u00d740=\"),\n    @(\"47\u00d766=\", \"68\u00d723=\"),\n    @(\"79\u00d783=\", \"83\u00d731=\"),\n    @(\"50\u00d769=\", \"82\u00d725=\"),\n    @(\"52\u00d759=\", \"40\u00d790=\")\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
